$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window geometry (best-effort; mirrors the maximized-window resize in the diff) ---
$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
$win.Width = 29040
$win.Height = 15840

# --- Header rename: "MSE" -> "Sredni blad [m]" (Polish: "Średni błąd [m]") ---
$ws.Range("G1").Value = "Średni błąd [m]"

# --- Row 3..13: replace MSE-style big numbers with the mean-distance-style small numbers ---
$ws.Range("G3").Value = 5.1212542330011503
$ws.Range("H3").Value = 5.70055632427436
$ws.Range("I3").Value = 2.8685542400222799

$ws.Range("G4").Value = 3.64492172581136
$ws.Range("H4").Value = 7.8288251408595002
$ws.Range("I4").Value = 2.8308090215864898

$ws.Range("G5").Value = 7.0710577745182297
$ws.Range("H5").Value = 7.57451590247853
$ws.Range("I5").Value = 3.4753383706745602

$ws.Range("G6").Value = 4.10133794441413
$ws.Range("H6").Value = 5.1066861988382399
$ws.Range("I6").Value = 2.7233334868241799

$ws.Range("G7").Value = 6.3848101617108197
$ws.Range("H7").Value = 5.1968136425635096
$ws.Range("I7").Value = 3.1215626569230701

$ws.Range("G8").Value = 4.58888371465234
$ws.Range("H8").Value = 5.7132084943813304
$ws.Range("I8").Value = 4.6267269125742603

$ws.Range("G9").Value = 4.2963607144986398
$ws.Range("H9").Value = 3.6083834453795101
$ws.Range("I9").Value = 2.5241108252828299

$ws.Range("G10").Value = 2.41677326657488
$ws.Range("H10").Value = 3.3237553242971298
$ws.Range("I10").Value = 3.3458391518149702

$ws.Range("G11").Value = 8.0934972050308005
$ws.Range("H11").Value = 6.7694905765076099
$ws.Range("I11").Value = 12.3995980610463

$ws.Range("G12").Value = 8.1189668947355607
$ws.Range("H12").Value = 11.2097337141926
$ws.Range("I12").Value = 9.5129165829411697

$ws.Range("G13").Value = 9.0785509484252493
$ws.Range("H13").Value = 10.472118526196899
$ws.Range("I13").Value = 10.0935328105118

# --- Selection moved to K26 ---
$ws.Range("K26").Select()
